$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: record that no one recovered / died that day (Recover / Death = 0) ---
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0

# --- New row 23: Kenyan Corona update for 5 April 2020 ---
$ws.Range("A23").Value = 43926
$ws.Range("A23").NumberFormat = "d-mmm-yy"
$ws.Range("A23").HorizontalAlignment = -4108

$ws.Range("B23").Value = 16
$ws.Range("C23").Value = 530
$ws.Range("D23").Value = "Nigerian(1)"
$ws.Range("E23").Value = "Nairobi(12), Kilifi(1),Mombasa(3)"
$ws.Range("F23").Value = 142
$ws.Range("K23").Value = "Mercy"
$ws.Range("L23").Value = "22-66"
$ws.Range("G23").Value = "Community(5), Imported(11)"
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("O23").Value = 7
$ws.Range("P23").Value = 9

# --- widen the "Travelled From" column slightly to fit the new entries ---
$ws.Columns("D").ColumnWidth = 26

# --- update the view: scroll back to A1 and move the active selection ---
[void]$ws.Range("I19").Select()
